# injection results - added average air pressure
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Ave Temp c (water)" row (row 32) to
# hold the new "Ave Pressure kpa (air)" measurement. Excel copies the
# formatting of the row above (row 31) into the freshly inserted row.
$ws.Rows("32:32").Insert()

$ws.Range("A32").Value = "Ave Pressure kpa (air)"
$ws.Range("B32").Value = 62.77

# The old highlight/fill formatting on A3:A5 / B3:B5 is no longer applied;
# clear the interior fill so those cells fall back to their plain styles.
$ws.Range("A3:A5").Interior.Pattern = -4142
$ws.Range("B3:B5").Interior.Pattern = -4142

# Restore the on-screen selection to match the edited area.
[void]$ws.Range("E33").Select()
